$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.694.23'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.599.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("E8").Value = '  +0.68%  '

$ws.Range("E9").Value = '  +1.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.50'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0842'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.823.80'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.599.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.50%  '

$ws.Range("E15").Value = '  +0.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.662.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.16%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0764'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '210.05'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.48%  '

$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.66%  '

$ws.Range("E22").Value = '  +1.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("E24").Value = '  +1.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.70%  '

$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.34%  '

$ws.Range("E28").Value = '  +0.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.85%  '

$ws.Range("E30").Value = '  +3.36%  '

$ws.Range("E31").Value = '  +0.29%  '

$ws.Range("E32").Value = '  +0.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.284.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.619'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.60%  '

$ws.Range("E37").Value = '  +0.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0172'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.27%  '

$ws.Range("E39").Value = '  +17.84%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.826'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.45'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.784'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.26%  '

$ws.Range("E43").Value = '  -0.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.735.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.40%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.36%  '

$ws.Range("E47").Value = '  -2.14%  '

$ws.Range("E48").Value = '  +0.13%  '

$ws.Range("E49").Value = '  +0.46%  '

$ws.Range("E50").Value = '  +0.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.60%  '
